# The "Company B" record (originally row 3: Company B / POC B / DESIGNATION B /
# dronekill1604@gmail.com) was removed from the tracker sheet. Deleting the
# whole row shifts the "Company C" record up from row 4 to row 3, and Excel
# automatically drops the now-unused shared strings for the removed record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("3").Delete()
